$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper pattern: date-like strings (dd.mm.yyyy / yyyy-mm-dd) must stay as
# plain text cells (matching the source inlineStr cells), not be
# auto-converted into Excel date serials. Forcing NumberFormat to Text
# before the assignment keeps the literal string, and resetting the
# Style back to "Normal" afterwards restores the original (unstyled)
# cell format so no stray style index is introduced.

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 8: advance dates by one day
Set-TextValue $ws.Range("A8") "11.08.2021"
Set-TextValue $ws.Range("B8") "2021-08-11"

# Row 13 and 14: swap the D/E/F/H content between the two rows
$ws.Range("D13").Value = "тест 11"
$ws.Range("E13").Value = "тест 11"
$ws.Range("F13").Value = "тест 11"
$ws.Range("H13").Value = 13

$ws.Range("D14").Value = "Шарометный бой"
$ws.Range("E14").Value = "Дарт Вейдер"
$ws.Range("F14").Value = "Уничтожение повстанцев"
$ws.Range("H14").Value = 19

# Row 17: rename title
$ws.Range("D17").Value = "Игра в домино"

# Row 20: shift dates by one day
Set-TextValue $ws.Range("A20") "16.09.2021"
Set-TextValue $ws.Range("B20") "2021-09-16"
Set-TextValue $ws.Range("C20") "20.09.2021"

# Row 21: new dates + swapped content (now matches old row 22 content)
Set-TextValue $ws.Range("A21") "17.09.2021"
Set-TextValue $ws.Range("B21") "2021-09-17"
Set-TextValue $ws.Range("C21") "18.09.2021"
$ws.Range("D21").Value = "Самая лучшая игра"
$ws.Range("E21").Value = "Самый лучший организатор"
$ws.Range("F21").Value = "Самый лучший тип игр"
$ws.Range("H21").Value = 24

# Row 22: new dates + swapped content (now matches old row 21 content)
Set-TextValue $ws.Range("A22") "20.09.2021"
Set-TextValue $ws.Range("B22") "2021-09-20"
Set-TextValue $ws.Range("C22") "23.09.2021"
$ws.Range("D22").Value = "ИСАФ"
$ws.Range("E22").Value = "HQ"
$ws.Range("F22").Value = "милсим"
$ws.Range("H22").Value = 17
